$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Find the "Date" column header on row 1.
$dateCol = 0
for ($col = 1; $col -le $lastCol; $col++) {
    if ($ws.Cells.Item(1, $col).Text -eq "Date") {
        $dateCol = $col
    }
}

if ($dateCol -gt 0) {
    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, $dateCol)
        if ($cell.Text -eq "5-18-2007-08") {
            # Assign through a literal-text formula + paste-as-value so the
            # result stays a plain text string (not auto-parsed as a date)
            # and keeps the cell's original (default) number format/style.
            $cell.Formula = '="2008-05-18"'
            $cell.Copy() | Out-Null
            $cell.PasteSpecial(-4163) | Out-Null
        }
    }
}
